$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.05441666666666667
$ws.Cells.Item(2, 8).Value = 0.16325
$ws.Cells.Item(2, 9).Value = 0.00608027172874025
$ws.Cells.Item(2, 10).Value = 0.006080271728740251
$ws.Cells.Item(2, 13).Value = 13.00733333333333
$ws.Cells.Item(2, 14).Value = 39.022
$ws.Cells.Item(2, 15).Value = 0.1070948256796854
$ws.Cells.Item(2, 16).Value = 0.1106092800596386
$ws.Cells.Item(2, 17).Value = 0.7078157222222222
$ws.Cells.Item(2, 18).Value = 6.3703415
$ws.Cells.Item(2, 19).Value = 0.0006511656408745563
$ws.Cells.Item(2, 20).Value = 0.0006725344784829336
$ws.Cells.Item(3, 7).Value = 0.05441666666666667
$ws.Cells.Item(3, 8).Value = 0.16325
$ws.Cells.Item(3, 9).Value = 0.00608027172874025
$ws.Cells.Item(3, 10).Value = 0.006080271728740251
$ws.Cells.Item(3, 15).Value = 0.7923813923330444
$ws.Cells.Item(3, 16).Value = 0.8183844063648099
$ws.Cells.Item(3, 17).Value = 5.237041135555557
$ws.Cells.Item(3, 18).Value = 47.13337022000001
$ws.Cells.Item(3, 19).Value = 0.004817894178182447
$ws.Cells.Item(3, 20).Value = 0.004975999569261827
$ws.Cells.Item(4, 7).Value = 0.05441666666666667
$ws.Cells.Item(4, 8).Value = 0.16325
$ws.Cells.Item(4, 9).Value = 0.00608027172874025
$ws.Cells.Item(4, 10).Value = 0.006080271728740251
$ws.Cells.Item(4, 13).Value = 0.298136
$ws.Cells.Item(4, 14).Value = 0.8944080000000001
$ws.Cells.Item(4, 15).Value = 0.002454678613257035
$ws.Cells.Item(4, 16).Value = 0.002535232047552183
$ws.Cells.Item(4, 17).Value = 0.01622356733333333
$ws.Cells.Item(4, 18).Value = 0.146012106
$ws.Cells.Item(4, 19).Value = 0.00001492511297533007
$ws.Cells.Item(4, 20).Value = 0.0000154148997445278
$ws.Cells.Item(5, 7).Value = 0.05441666666666667
$ws.Cells.Item(5, 8).Value = 0.16325
$ws.Cells.Item(5, 9).Value = 0.00608027172874025
$ws.Cells.Item(5, 10).Value = 0.006080271728740251
$ws.Cells.Item(5, 13).Value = 11.5773025
$ws.Cells.Item(5, 14).Value = 23.154605
$ws.Cells.Item(5, 15).Value = 0.0953207826158438
$ws.Cells.Item(5, 16).Value = 0.06563257109105912
$ws.Cells.Item(5, 17).Value = 0.6299982110416666
$ws.Cells.Item(5, 18).Value = 3.77998926625
$ws.Cells.Item(5, 19).Value = 0.0005795762597005102
$ws.Cells.Item(5, 20).Value = 0.0003990638664895015
$ws.Cells.Item(6, 7).Value = 0.05441666666666667
$ws.Cells.Item(6, 8).Value = 0.16325
$ws.Cells.Item(6, 9).Value = 0.00608027172874025
$ws.Cells.Item(6, 10).Value = 0.006080271728740251
$ws.Cells.Item(6, 13).Value = 0.3338006666666667
$ws.Cells.Item(6, 14).Value = 1.001402
$ws.Cells.Item(6, 15).Value = 0.00274832075816945
$ws.Cells.Item(6, 16).Value = 0.002838510436940246
$ws.Cells.Item(6, 17).Value = 0.01816431961111111
$ws.Cells.Item(6, 18).Value = 0.1634788765
$ws.Cells.Item(6, 19).Value = 0.00001671053700740768
$ws.Cells.Item(6, 20).Value = 0.00001725891476146191
$ws.Cells.Item(7, 7).Value = 4.046611333333334
$ws.Cells.Item(7, 8).Value = 12.139834
$ws.Cells.Item(7, 9).Value = 0.4521500120171497
$ws.Cells.Item(7, 10).Value = 0.4521500120171497
$ws.Cells.Item(7, 13).Value = 13.00733333333333
$ws.Cells.Item(7, 14).Value = 39.022
$ws.Cells.Item(7, 15).Value = 0.1070948256796854
$ws.Cells.Item(7, 16).Value = 0.1106092800596386
$ws.Cells.Item(7, 17).Value = 52.63562248311112
$ws.Cells.Item(7, 18).Value = 473.720602348
$ws.Cells.Item(7, 19).Value = 0.04842292671804429
$ws.Cells.Item(7, 20).Value = 0.05001198730817388
$ws.Cells.Item(8, 7).Value = 4.046611333333334
$ws.Cells.Item(8, 8).Value = 12.139834
$ws.Cells.Item(8, 9).Value = 0.4521500120171497
$ws.Cells.Item(8, 10).Value = 0.4521500120171497
$ws.Cells.Item(8, 15).Value = 0.7923813923330444
$ws.Cells.Item(8, 16).Value = 0.8183844063648099
$ws.Cells.Item(8, 17).Value = 389.4444718947379
$ws.Cells.Item(8, 18).Value = 3505.000247052641
$ws.Cells.Item(8, 19).Value = 0.3582752560655518
$ws.Cells.Item(8, 20).Value = 0.3700325191724967
$ws.Cells.Item(9, 7).Value = 4.046611333333334
$ws.Cells.Item(9, 8).Value = 12.139834
$ws.Cells.Item(9, 9).Value = 0.4521500120171497
$ws.Cells.Item(9, 10).Value = 0.4521500120171497
$ws.Cells.Item(9, 13).Value = 0.298136
$ws.Cells.Item(9, 14).Value = 0.8944080000000001
$ws.Cells.Item(9, 15).Value = 0.002454678613257035
$ws.Cells.Item(9, 16).Value = 0.002535232047552183
$ws.Cells.Item(9, 17).Value = 1.206440516474667
$ws.Cells.Item(9, 18).Value = 10.857964648272
$ws.Cells.Item(9, 19).Value = 0.001109882964482409
$ws.Cells.Item(9, 20).Value = 0.001146305200766982
$ws.Cells.Item(10, 7).Value = 4.046611333333334
$ws.Cells.Item(10, 8).Value = 12.139834
$ws.Cells.Item(10, 9).Value = 0.4521500120171497
$ws.Cells.Item(10, 10).Value = 0.4521500120171497
$ws.Cells.Item(10, 13).Value = 11.5773025
$ws.Cells.Item(10, 14).Value = 23.154605
$ws.Cells.Item(10, 15).Value = 0.0953207826158438
$ws.Cells.Item(10, 16).Value = 0.06563257109105912
$ws.Cells.Item(10, 17).Value = 46.84884350592834
$ws.Cells.Item(10, 18).Value = 281.09306103557
$ws.Cells.Item(10, 19).Value = 0.04309929300523788
$ws.Cells.Item(10, 20).Value = 0.02967576780753881
$ws.Cells.Item(11, 7).Value = 4.046611333333334
$ws.Cells.Item(11, 8).Value = 12.139834
$ws.Cells.Item(11, 9).Value = 0.4521500120171497
$ws.Cells.Item(11, 10).Value = 0.4521500120171497
$ws.Cells.Item(11, 13).Value = 0.3338006666666667
$ws.Cells.Item(11, 14).Value = 1.001402
$ws.Cells.Item(11, 15).Value = 0.00274832075816945
$ws.Cells.Item(11, 16).Value = 0.002838510436940246
$ws.Cells.Item(11, 17).Value = 1.350761560807556
$ws.Cells.Item(11, 18).Value = 12.156854047268
$ws.Cells.Item(11, 19).Value = 0.001242653263833299
$ws.Cells.Item(11, 20).Value = 0.001283432528173337
$ws.Cells.Item(12, 7).Value = 4.848681666666667
$ws.Cells.Item(12, 8).Value = 14.546045
$ws.Cells.Item(12, 9).Value = 0.5417697162541101
$ws.Cells.Item(12, 10).Value = 0.5417697162541102
$ws.Cells.Item(12, 13).Value = 13.00733333333333
$ws.Cells.Item(12, 14).Value = 39.022
$ws.Cells.Item(12, 15).Value = 0.1070948256796854
$ws.Cells.Item(12, 16).Value = 0.1106092800596386
$ws.Cells.Item(12, 17).Value = 63.06841866555555
$ws.Cells.Item(12, 18).Value = 567.61576799
$ws.Cells.Item(12, 19).Value = 0.05802073332076652
$ws.Cells.Item(12, 20).Value = 0.05992475827298183
$ws.Cells.Item(13, 7).Value = 4.848681666666667
$ws.Cells.Item(13, 8).Value = 14.546045
$ws.Cells.Item(13, 9).Value = 0.5417697162541101
$ws.Cells.Item(13, 10).Value = 0.5417697162541102
$ws.Cells.Item(13, 15).Value = 0.7923813923330444
$ws.Cells.Item(13, 16).Value = 0.8183844063648099
$ws.Cells.Item(13, 17).Value = 466.635442723689
$ws.Cells.Item(13, 18).Value = 4199.718984513201
$ws.Cells.Item(13, 19).Value = 0.4292882420893102
$ws.Cells.Item(13, 20).Value = 0.4433758876230515
$ws.Cells.Item(14, 7).Value = 4.848681666666667
$ws.Cells.Item(14, 8).Value = 14.546045
$ws.Cells.Item(14, 9).Value = 0.5417697162541101
$ws.Cells.Item(14, 10).Value = 0.5417697162541102
$ws.Cells.Item(14, 13).Value = 0.298136
$ws.Cells.Item(14, 14).Value = 0.8944080000000001
$ws.Cells.Item(14, 15).Value = 0.002454678613257035
$ws.Cells.Item(14, 16).Value = 0.002535232047552183
$ws.Cells.Item(14, 17).Value = 1.445566557373333
$ws.Cells.Item(14, 18).Value = 13.01009901636
$ws.Cells.Item(14, 19).Value = 0.001329870535799297
$ws.Cells.Item(14, 20).Value = 0.001373511947040673
$ws.Cells.Item(15, 7).Value = 4.848681666666667
$ws.Cells.Item(15, 8).Value = 14.546045
$ws.Cells.Item(15, 9).Value = 0.5417697162541101
$ws.Cells.Item(15, 10).Value = 0.5417697162541102
$ws.Cells.Item(15, 13).Value = 11.5773025
$ws.Cells.Item(15, 14).Value = 23.154605
$ws.Cells.Item(15, 15).Value = 0.0953207826158438
$ws.Cells.Item(15, 16).Value = 0.06563257109105912
$ws.Cells.Item(15, 17).Value = 56.13465438120417
$ws.Cells.Item(15, 18).Value = 336.807926287225
$ws.Cells.Item(15, 19).Value = 0.05164191335090541
$ws.Cells.Item(15, 20).Value = 0.03555773941703082
$ws.Cells.Item(16, 7).Value = 4.848681666666667
$ws.Cells.Item(16, 8).Value = 14.546045
$ws.Cells.Item(16, 9).Value = 0.5417697162541101
$ws.Cells.Item(16, 10).Value = 0.5417697162541102
$ws.Cells.Item(16, 13).Value = 0.3338006666666667
$ws.Cells.Item(16, 14).Value = 1.001402
$ws.Cells.Item(16, 15).Value = 0.00274832075816945
$ws.Cells.Item(16, 16).Value = 0.002838510436940246
$ws.Cells.Item(16, 17).Value = 1.618493172787778
$ws.Cells.Item(16, 18).Value = 14.56643855509
$ws.Cells.Item(16, 19).Value = 0.001488956957328744
$ws.Cells.Item(16, 20).Value = 0.001725891476146191
